$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'72.680.99"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +2.85%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'3.980.83"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +1.42%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.24%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'587.20"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +8.66%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'157.69"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +6.71%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.681"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.52%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.749"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +1.95%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'0.169"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +1.33%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'53.28"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.60%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'0.0000322"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +1.38%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'10.87"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +3.70%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'4.620.25"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +1.24%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'3.969.13"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.09%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  +9.52%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'14.08"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +1.31%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'20.42"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +1.00%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +0.02%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'72.448.94"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +2.69%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'433.64"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +2.02%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'4.71"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +11.52%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'96.06"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.20%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'3.43"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -2.21%  "

# Row 25
$ws.Cells.Item(25, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(25, 4).Value = "'14.35"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.28%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "Toncoin"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(26, 4).Value = "'4.44"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +22.95%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'11.13"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.68%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'10.74"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +2.49%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'5.93"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +1.47%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'36.50"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.95%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'7.82"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +6.05%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'13.58"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.88%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +1.95%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'48.82"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +2.07%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'680.32"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.54%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'68.25"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +5.00%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "PEPE"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(37, 4).Value = "'0.0₃0887"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +8.82%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "TheGraph"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(38, 4).Value = "'0.437"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +1.66%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "WEMIXToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(39, 4).Value = "'3.36"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.64%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "ThetaToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Cells.Item(40, 4).Value = "'3.36"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.06%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "Kaspa"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(41, 4).Value = "'0.146"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -1.57%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "Dai"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(42, 4).Value = "'1.00"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.03%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -0.08%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "VeChain"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(44, 4).Value = "'0.0486"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +1.54%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "THORChain"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Cells.Item(45, 4).Value = "'10.78"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +12.56%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'0.149"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.79%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'2.66"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -1.20%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'3.39"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +0.54%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'3.03"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +2.48%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'3.43"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +6.16%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  +9.04%  "
